$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J, matching the style used by existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-21
$data = @(
    @(1, 6),
    @(1, 4),
    @(5, 8),
    @(1, 5),
    @(1, 2),
    @(1, 3),
    @(1, 2),
    @(1, 4),
    @(1, 6),
    @(1, 6),
    @(7, 7),
    @(7, 7),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(6, 8),
    @(6, 8),
    @(5, 7),
    @(3, 4)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
